$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Coin/Link/Price/Volume columns keep their original text
# representation (e.g. "30.20", "1.965.75") instead of being
# auto-converted to numbers by Excel when the value looks numeric.
$ws.Range("B2:E51").NumberFormat = "@"

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '31.843.74', '  +6.74%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.720.94', '  +5.01%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.998', '  -0.23%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '223.76', '  +3.91%  '),
    @(6, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.541', '  +4.12%  '),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.997', '  -0.24%  '),
    @(8, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '30.20', '  +4.66%  '),
    @(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '45.36', '  +3.59%  '),
    @(10, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.271', '  +4.20%  '),
    @(11, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.0654', '  +7.49%  '),
    @(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.0911', '  +1.35%  '),
    @(13, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.965.75', '  +4.95%  '),
    @(14, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.715.40', '  +4.69%  '),
    @(15, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.618', '  +4.26%  '),
    @(16, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '10.22', '  +6.69%  '),
    @(17, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.21', '  +7.86%  '),
    @(18, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '31.789.62', '  +6.57%  '),
    @(19, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '67.63', '  +5.18%  '),
    @(20, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '253.68', '  +6.74%  '),
    @(21, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0₃0727', '  +3.40%  '),
    @(22, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.998', '  -0.21%  '),
    @(23, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '10.18', '  +2.66%  '),
    @(24, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.28', '  +3.49%  '),
    @(25, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.16', '  -0.90%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '159.38', '  +1.18%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '16.22', '  +4.05%  '),
    @(28, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.113', '  +3.51%  '),
    @(29, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '6.86', '  +3.19%  '),
    @(30, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.998', '  -0.19%  '),
    @(31, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.87', '  +14.49%  '),
    @(32, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0506', '  +2.16%  '),
    @(33, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.17', '  +5.15%  '),
    @(34, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.43', '  +7.55%  '),
    @(35, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.528.17', '  +7.60%  '),
    @(36, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.76', '  +4.48%  '),
    @(37, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.05', '  +3.07%  '),
    @(38, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '83.48', '  +8.96%  '),
    @(39, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.613', '  +8.33%  '),
    @(40, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0182', '  +5.02%  '),
    @(41, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.73', '  +2.50%  '),
    @(42, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.31', '  +0.72%  '),
    @(43, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '0.857', '  +2.94%  '),
    @(44, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.04', '  +5.25%  '),
    @(45, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.0503', '  +0.06%  '),
    @(46, 'BitcoinSV', 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv', '54.10', '  +7.74%  '),
    @(47, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.03', '  +3.51%  '),
    @(48, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.997', '  -0.25%  '),
    @(49, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.61', '  +5.27%  '),
    @(50, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '1.858.99', '  +4.39%  '),
    @(51, 'BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.0₆0119', '  +7.32%  '),
)

foreach ($item in $data) {
    $row = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
}
